$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing row 5 (previously the "green block" question) ---
# A5 text changes to the "Canadian flag" colour phrasing, goal values shift, category becomes 2
$ws.Range("A5").Value = "Take the block with a colour on the Canadian flag, and move it to the position (0,0.375). "
$ws.Range("B5").Value = "(0,0.375,0.02)"
$ws.Range("C5").Value = "(-0.25,0.25,0.02)"
$ws.Range("D5").Value = "(0.25,0.5,0.02)"
$ws.Range("E5").Value = "(-0.25,0.5,0.02)"
$ws.Range("F5").Value = 2

# --- New data rows 6-19 ---
$rows = @(
  @{ r=6;  a="Take the block that is the same colour as the ocean, and move it to the position (0,0.375). "; b="(0.25,0.25,0.02)"; c="(0,0.375,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=2 },
  @{ r=7;  a="Take the block that is the colour of the Saharah Desert, and move it to the position (0,0.375). "; b="(0.25,0.25,0.02)"; c="(-0.25,0.25,0.02)"; d="(0,0.375,0.02)"; e="(-0.25,0.5,0.02)"; f=2 },
  @{ r=8;  a="Take the blue block and move it directly in between the red and yellow blocks"; b="(0.25,0.25,0.02)"; c="(0,0.5,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=3 },
  @{ r=9;  a="Take the red block and move it directly in between the green and yellow blocks"; b="(0,0.5,0.02)"; c="(-0.25,0.25,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=3 },
  @{ r=10; a="Take the yellow block and move it directly between the green and blue blocks"; b="(0.25,0.25,0.02)"; c="(-0.25,0.25,0.02)"; d="(-0.25,0.375,0.02)"; e="(-0.25,0.5,0.02)"; f=3 },
  @{ r=11; a="Pick up the blue block and move it directly in between the green and yellow blocks"; b="(0.25,0.25,0.02)"; c="(0,0.5,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=4 },
  @{ r=12; a="Pick up the red block and move it directly in between the green and yellow blocks"; b="(0,0.5,0.02)"; c="(-0.25,0.25,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=4 },
  @{ r=13; a="Pick up the yellow block and move it directly between the green and blue blocks"; b="(0.25,0.25,0.02)"; c="(-0.25,0.25,0.02)"; d="(-0.25,0.375,0.02)"; e="(-0.25,0.5,0.02)"; f=4 },
  @{ r=14; a="Pick up the ocean-colored block and move it directly in between the green and yellow blocks"; b="(0.25,0.25,0.02)"; c="(0,0.5,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=5 },
  @{ r=15; a="Pick up the block that is the color of blood and move it directly in between the green and yellow blocks"; b="(0,0.5,0.02)"; c="(-0.25,0.25,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=5 },
  @{ r=16; a="Pick up the block that is the color of a lemon and move it directly between the green and blue blocks"; b="(0.25,0.25,0.02)"; c="(-0.25,0.25,0.02)"; d="(-0.25,0.375,0.02)"; e="(-0.25,0.5,0.02)"; f=5 },
  @{ r=17; a="Pick up the block that has the color of the sky and move it directly in between the forest colored block and the lemon colored block"; b="(0.25,0.25,0.02)"; c="(0,0.5,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=6 },
  @{ r=18; a="Pick up the block that is the color of a stop sign and move it directly in between the grass colored block and block with color closest to a banana"; b="(0,0.5,0.02)"; c="(-0.25,0.25,0.02)"; d="(0.25,0.5,0.02)"; e="(-0.25,0.5,0.02)"; f=6 },
  @{ r=19; a="Pick up the block that is the color of a lemon and move it directly between the block with the color of grass and the block with the color of the sky"; b="(0.25,0.25,0.02)"; c="(-0.25,0.25,0.02)"; d="(-0.25,0.375,0.02)"; e="(-0.25,0.5,0.02)"; f=6 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.a
    $ws.Range("B$r").Value = $row.b
    $ws.Range("C$r").Value = $row.c
    $ws.Range("D$r").Value = $row.d
    $ws.Range("E$r").Value = $row.e
    $ws.Range("F$r").Value = $row.f
}

# --- Update the selection to match the post-edit workbook state ---
$ws.Range("A22").Select()
